$wb = $excel.ActiveWorkbook

# --- Sheet: "Add Devices" -------------------------------------------------
$wsAdd = $wb.Worksheets.Item("Add Devices")
$wsAdd.Range("L10").Value = $false
$wsAdd.Range("M10").Value = $false
$wsAdd.Range("M11").Value = $false

# --- Sheet: "Delete Devices" ----------------------------------------------
$wsDel = $wb.Worksheets.Item("Delete Devices")
$wsDel.Range("G8").Value = $false
$wsDel.Range("H8").Value = $false

# --- Sheet: "Add Devices for Second Panel" --------------------------------
$wsAdd2 = $wb.Worksheets.Item("Add Devices for Second Panel")
$wsAdd2.Range("L10").Value = $false
$wsAdd2.Range("M10").Value = $false
$wsAdd2.Range("M11").Value = $false

# --- Sheet: "Delete Devices for Second Panel" -----------------------------
$wsDel2 = $wb.Worksheets.Item("Delete Devices for Second Panel")
$wsDel2.Range("G8").Value = $false
$wsDel2.Range("H8").Value = $false

# Update per-sheet selections to match what was left behind when the
# workbook was saved.
$wsDel.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$wsDel.Range("H8").Select() | Out-Null

$wsAdd2.Range("M11").Select() | Out-Null

$wsDel2.Range("H8").Select() | Out-Null

# "Add Devices" is the sheet that ends up active/selected in the saved
# workbook, so activate it last (moves tabSelected/activeTab onto it and
# off "Delete Devices for Second Panel").
$wsAdd.Activate() | Out-Null
$wsAdd.Range("M11").Select() | Out-Null
